$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A60").Value = "Davide Raffaelli "
$ws.Range("B60").Value = "Alberto Cerisara | SHARK ATTACK"
$ws.Range("C60").Value = "Daniele Feller | GREP"
$ws.Range("D60").Value = "Matteo Simoncelli | IMONTAGNA"
$ws.Range("E60").Value = "Leonardo Viola | SHARK ATTACK"
$ws.Range("F60").Value = "Alessio Debiasi | Mai una gioia"
